# Refresh Leve profit-calculation columns (currentAveragePrice(NQ/HQ),
# LevePrice(NQ/HQ), LeveProfit(NQ/HQ)) with the latest Universalis market
# data snapshot for the affected Leve rows across all job sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 360.66666
$ws.Range("I12").Value = 425
$ws.Range("K12").Value = 425
$ws.Range("M12").Value = -255
$ws.Range("H28").Value = 627.2353000000001
$ws.Range("I28").Value = 577.26666
$ws.Range("K28").Value = 577.26666
$ws.Range("M28").Value = -92.26666
$ws.Range("H33").Value = 501.5
$ws.Range("I33").Value = 795
$ws.Range("K33").Value = 795
$ws.Range("M33").Value = -566
$ws.Range("H39").Value = 1212.2142
$ws.Range("I39").Value = 406.45456
$ws.Range("J39").Value = 4166.6665
$ws.Range("K39").Value = 1219.36368
$ws.Range("L39").Value = 12499.9995
$ws.Range("M39").Value = -923.3636799999999
$ws.Range("N39").Value = -13091.9995
$ws.Range("H52").Value = 1496
$ws.Range("I52").Value = 1607.1428
$ws.Range("K52").Value = 4821.428400000001
$ws.Range("M52").Value = -4661.428400000001
$ws.Range("H104").Value = 6945.1665
$ws.Range("I104").Value = 6945.1665
$ws.Range("K104").Value = 20835.4995
$ws.Range("M104").Value = -19088.4995
$ws.Range("H137").Value = 3272.8
$ws.Range("I137").Value = 2962.375
$ws.Range("K137").Value = 8887.125
$ws.Range("M137").Value = -6337.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 34157
$ws.Range("I28").Value = 34157
$ws.Range("K28").Value = 34157
$ws.Range("M28").Value = -33965
$ws.Range("H32").Value = 35641.61
$ws.Range("I32").Value = 37063.758
$ws.Range("K32").Value = 37063.758
$ws.Range("M32").Value = -36776.758
$ws.Range("H61").Value = 12416.8
$ws.Range("I61").Value = 11879.066
$ws.Range("J61").Value = 14030
$ws.Range("K61").Value = 11879.066
$ws.Range("L61").Value = 14030
$ws.Range("M61").Value = -11667.066
$ws.Range("N61").Value = -14454
$ws.Range("H74").Value = 1004970.9
$ws.Range("I74").Value = 1669668.1
$ws.Range("J74").Value = 7925
$ws.Range("K74").Value = 1669668.1
$ws.Range("L74").Value = 7925
$ws.Range("M74").Value = -1668794.1
$ws.Range("N74").Value = -9673
$ws.Range("H77").Value = 1004970.9
$ws.Range("I77").Value = 1669668.1
$ws.Range("J77").Value = 7925
$ws.Range("K77").Value = 8348340.5
$ws.Range("L77").Value = 39625
$ws.Range("M77").Value = -8343972.5
$ws.Range("N77").Value = -48361
$ws.Range("H99").Value = 34157
$ws.Range("I99").Value = 34157
$ws.Range("K99").Value = 34157
$ws.Range("M99").Value = -31162
$ws.Range("H136").Value = 12416.8
$ws.Range("I136").Value = 11879.066
$ws.Range("J136").Value = 14030
$ws.Range("K136").Value = 35637.198
$ws.Range("L136").Value = 42090
$ws.Range("M136").Value = -33087.198
$ws.Range("N136").Value = -47190

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 3426
$ws.Range("I22").Value = 3900.6667
$ws.Range("K22").Value = 3900.6667
$ws.Range("M22").Value = -3727.6667
$ws.Range("H64").Value = 3663.3333
$ws.Range("J64").Value = 3663.3333
$ws.Range("L64").Value = 3663.3333
$ws.Range("N64").Value = -4113.3333
$ws.Range("H67").Value = 3663.3333
$ws.Range("J67").Value = 3663.3333
$ws.Range("L67").Value = 3663.3333
$ws.Range("N67").Value = -5223.3333
$ws.Range("H86").Value = 2153.5
$ws.Range("I86").Value = 2153.5
$ws.Range("K86").Value = 2153.5
$ws.Range("M86").Value = -1030.5
$ws.Range("H89").Value = 2153.5
$ws.Range("I89").Value = 2153.5
$ws.Range("K89").Value = 10767.5
$ws.Range("M89").Value = -5151.5
$ws.Range("H107").Value = 1151.6538
$ws.Range("J107").Value = 1271.5454
$ws.Range("L107").Value = 1271.5454
$ws.Range("N107").Value = -5111.5454
$ws.Range("H134").Value = 4658.3257
$ws.Range("I134").Value = 3568.3333
$ws.Range("J134").Value = 8255.299999999999
$ws.Range("K134").Value = 10704.9999
$ws.Range("L134").Value = 24765.9
$ws.Range("M134").Value = -8169.999899999999
$ws.Range("N134").Value = -29835.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()
$ws.Range("H58").Value = 4159.3438
$ws.Range("I58").Value = 2738.25
$ws.Range("J58").Value = 5580.4375
$ws.Range("K58").Value = 2738.25
$ws.Range("L58").Value = 5580.4375
$ws.Range("M58").Value = -2535.25
$ws.Range("N58").Value = -5986.4375
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H132").Value = 5227.4614
$ws.Range("I132").Value = 3493.8572
$ws.Range("J132").Value = 7250
$ws.Range("K132").Value = 10481.5716
$ws.Range("L132").Value = 21750
$ws.Range("M132").Value = -7951.571599999999
$ws.Range("N132").Value = -26810
$ws.Range("H134").Value = 3730.8333
$ws.Range("I134").Value = 2747.2727
$ws.Range("K134").Value = 8241.8181
$ws.Range("M134").Value = -5706.8181
$ws.Range("H136").Value = 4159.3438
$ws.Range("I136").Value = 2738.25
$ws.Range("J136").Value = 5580.4375
$ws.Range("K136").Value = 8214.75
$ws.Range("L136").Value = 16741.3125
$ws.Range("M136").Value = -5664.75
$ws.Range("N136").Value = -21841.3125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 1842.7778
$ws.Range("I40").Value = 275
$ws.Range("J40").Value = 2626.6667
$ws.Range("K40").Value = 1100
$ws.Range("L40").Value = 10506.6668
$ws.Range("M40").Value = -1031
$ws.Range("N40").Value = -10644.6668
$ws.Range("H80").Value = 3679.9
$ws.Range("J80").Value = 3811
$ws.Range("L80").Value = 11433
$ws.Range("N80").Value = -13305
$ws.Range("H83").Value = 3679.9
$ws.Range("J83").Value = 3811
$ws.Range("L83").Value = 34299
$ws.Range("N83").Value = -43659
$ws.Range("H141").Value = 8325.833000000001
$ws.Range("I141").Value = 8597.25
$ws.Range("K141").Value = 25791.75
$ws.Range("M141").Value = -20611.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 69893.28999999999
$ws.Range("J57").Value = 69893.28999999999
$ws.Range("L57").Value = 69893.28999999999
$ws.Range("N57").Value = -71533.28999999999
$ws.Range("H70").Value = 5829
$ws.Range("I70").Value = 5004
$ws.Range("K70").Value = 5004
$ws.Range("M70").Value = -4734
$ws.Range("H73").Value = 5829
$ws.Range("I73").Value = 5004
$ws.Range("K73").Value = 5004
$ws.Range("M73").Value = -4068
$ws.Range("H104").Value = 119999.5
$ws.Range("J104").Value = 119999.5
$ws.Range("L104").Value = 119999.5
$ws.Range("N104").Value = -126987.5
$ws.Range("H132").Value = 8496.375
$ws.Range("I132").Value = 6230
$ws.Range("J132").Value = 11410.286
$ws.Range("K132").Value = 18690
$ws.Range("L132").Value = 34230.858
$ws.Range("M132").Value = -16160
$ws.Range("N132").Value = -39290.858

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1279
$ws.Range("I16").Value = 1198.8889
$ws.Range("K16").Value = 1198.8889
$ws.Range("M16").Value = -1028.8889
$ws.Range("H55").Value = 145
$ws.Range("I55").Value = 121.5
$ws.Range("K55").Value = 121.5
$ws.Range("M55").Value = 51.5
$ws.Range("H140").Value = 85749.75
$ws.Range("J140").Value = 105999.664
$ws.Range("L140").Value = 105999.664
$ws.Range("N140").Value = -116359.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 11804.907
$ws.Range("I81").Value = 3063.4707
$ws.Range("K81").Value = 6126.9414
$ws.Range("M81").Value = -5065.9414
$ws.Range("H84").Value = 11804.907
$ws.Range("I84").Value = 3063.4707
$ws.Range("K84").Value = 30634.707
$ws.Range("M84").Value = -25330.707
$ws.Range("H113").Value = 8998.333000000001
$ws.Range("J113").Value = 8998
$ws.Range("L113").Value = 26994
$ws.Range("N113").Value = -31334
$ws.Range("H122").Value = 2605.5557
$ws.Range("I122").Value = 2556.3125
$ws.Range("K122").Value = 7668.9375
$ws.Range("M122").Value = -5218.9375
